# Update the "Metadata" sheet (StructureDefinition publishing info).
$wb = $excel.ActiveWorkbook
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/episode-disease-stage-code"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

# Clear the stale "Constraints" text that had leaked onto the root
# "Extension" row (AI2) on the "Elements" sheet - it belongs only to the
# "Extension.extension" row (AI4).
$elements = $wb.Worksheets.Item("Elements")
$elements.Range("AI2").Value = ""
